$wb = $excel.ActiveWorkbook

# Duplicate the "findCarTest" sheet (same brand/browser data) to create the
# new "findCarModelAndPriceTest" sheet, placing it right after the source.
$srcSheet = $wb.Worksheets.Item("findCarTest")
$srcSheet.Copy($null, $srcSheet)

# The copy gets auto-named "findCarTest (2)" and becomes the active sheet -
# rename it and set its own selection.
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "findCarModelAndPriceTest"
$newSheet.Range("F18").Select()

# Restore the originating sheet's selection (no longer the active tab).
$srcSheet.Range("A1:B5").Select()

$newSheet.Activate()
